# Mark additional completed tasks (checkmarks) for several students, per the
# "Download Feature & File Update" commit: cells that were previously blank
# in the "Tugas" (task) columns now get the Wingdings checkmark glyph "ü"
# (same glyph/style used by all the other "done" cells in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Checkmark($addr) {
    $rng = $ws.Range($addr)
    $rng.Value = "ü"
    $rng.Font.Name = "Wingdings"
    $rng.Font.Size = 12
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4108     # xlCenter
    $rng.Borders.LineStyle = 1         # xlContinuous
    $rng.Borders.Weight = 2            # xlThin
}

Set-Checkmark "E10"
Set-Checkmark "E12"
Set-Checkmark "E17"
Set-Checkmark "C22"
Set-Checkmark "E22"
Set-Checkmark "E29"

Write-Output "Checkmarks applied to E10, E12, E17, C22, E22, E29"
